# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to the way NBA stats were shown.
# Column BF ("Date") rows 2-31 were stored as "6-3-2012-13" and should
# read "2013-06-03" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("BF2:BF31")

# Make sure the values stay plain text (the source data is a text label
# that merely looks like a date) instead of being auto-converted into a
# real date serial number by Excel's input parser.
$range.NumberFormat = "@"

For ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 58).Value = "2013-06-03"
}
